$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-22 (Name, Soni, Narxi, Telefon raqami)
$data = @(
    @("Sadikov Anvar", 1, 320000, "+998946045808"),
    @("Muhammad", 1, 180000, "+998903203636"),
    @("Асадбек", 1, 320000, "+998970353930"),
    @("Muhammadali", 1, 300000, "+998977538191"),
    @("Малика Бадридинова", 1, 320000, "+998933833036"),
    @("Nozim", 1, 320000, "+998990833333"),
    @("Азиза", 1, 1200000, "+998900655055"),
    @("Шахзод", 1, 120000, "+998981223322"),
    @("Мухаммадамин", 1, 1100000, "+998979979191"),
    @("Мухлиса", 1, 150000, "+998981602444"),
    @("Миромон", 1, 380000, "+998998351771"),
    @("Абдурашид", 1, 320000, "+998957704244"),
    @("Бегзод ", 1, 750000, "+998977760464"),
    @("Sultanova Muborak", 1, 150000, "+998993656050"),
    @("Gulirano", 1, 120000, "+998958185775"),
    @("Raximov Muhammadjon", 1, 320000, "+998990390088"),
    @("Наима", 1, 400000, "+998977283121"),
    @("Амаль", 1, 320000, "+998959009050"),
    @("Nodir Qahramonov", 1, 300000, "+998915383553"),
    @("виталий", 1, 300000, "+998949446754"),
    @("Комрон", 1, 300000, "+998992132213")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]

    # Phone numbers start with "+" followed by digits, which Excel would
    # otherwise auto-convert to a number. Force text formatting, assign
    # the value, then restore the default "Normal" style so no stray
    # cell formatting is left behind.
    $phoneCell = $ws.Cells.Item($row, 4)
    $phoneCell.NumberFormat = "@"
    $phoneCell.Value = $data[$i][3]
    $phoneCell.Style = "Normal"
}

# Remove old rows 23-37 that are no longer present in the data
$lastOldRow = 37
$firstRemoveRow = $startRow + $data.Count
if ($lastOldRow -ge $firstRemoveRow) {
    $clearRange = $ws.Range("A$firstRemoveRow`:D$lastOldRow")
    $clearRange.ClearContents()
}
